# Reformat chord names in column B from the internal "Root[s]_suffix" naming
# scheme to a user-friendly display form, e.g. "Fs_min" -> "F#m", "A_maj" -> "A".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffixMap = @{
    'maj'  = ''
    'min'  = 'm'
    '7'    = '7'
    'min7' = 'm7'
}

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 414) { $lastRow = 414 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $raw = $cell.Text
    if ([string]::IsNullOrEmpty($raw)) { continue }

    if ($raw -match '^([A-G])(s)?_(maj|min7|min|7)$') {
        $root = $matches[1]
        $sharp = $matches[2]
        $suffix = $matches[3]

        $newRoot = $root
        if ($sharp -eq 's') { $newRoot = $root + '#' }

        $newVal = $newRoot + $suffixMap[$suffix]
        $cell.Value = $newVal
    }
}
